# Locate the page-break paragraph (the one whose text begins with a form-feed char \f, code 12)
$d = $word.ActiveDocument
$pageBreakParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.Length -gt 0 -and [int][char]$t[0] -eq 12) {
        $pageBreakParaIndex = $i
        break
    }
}
if ($pageBreakParaIndex -eq -1) {
    throw "Could not locate page-break paragraph"
}

$insertPos = $d.Paragraphs($pageBreakParaIndex).Range.Start
$insertRange = $d.Range($insertPos, $insertPos)

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
      <w:r>
        <w:t xml:space="preserve">Infine abbiamo </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>preso  qualche</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> misura a valori costanti per stimare l’errore</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">DC </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>staircase</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> ù</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Si è </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>fattio</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> variare il </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>th</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> da 0 a -30 </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>mV</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> per ogni valore di </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vbias</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> siamo partiti da 53,5V perché più sotto saremmo sotto </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Vbr</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>, e con passo di 0,5V siamo arrivati fino a 56V.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>ADC factor</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:strike/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:strike/>
        </w:rPr>
        <w:t>Psau</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:strike/>
        </w:rPr>
        <w:t xml:space="preserve"> a </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:strike/>
        </w:rPr>
        <w:t>th</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:strike/>
        </w:rPr>
        <w:t xml:space="preserve"> -5 </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:strike/>
        </w:rPr>
        <w:t>mV</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:strike/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:strike/>
        </w:rPr>
        <w:t>con  G</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:strike/>
        </w:rPr>
        <w:t xml:space="preserve"> = 10 non becca assolutamente un cazzo</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:strike/>
        </w:rPr>
        <w:t xml:space="preserve">, idem a </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:strike/>
        </w:rPr>
        <w:t>th</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:strike/>
        </w:rPr>
        <w:t xml:space="preserve"> -50, </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:strike/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:strike/>
        </w:rPr>
        <w:t xml:space="preserve">laser a 3 troppo basso, messo a </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Abbiamo messo il laser a 2 per non saturare il </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>sipm</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> e proceduto a trovare la prima configurazione misurabile</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> di gain e </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>thr</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
</w:p>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertRange.InsertXML($xml) | Out-Null

# Re-attach the "_GoBack" bookmark between "...misurabile" and " di gain e ..."
$findRange = $d.Content
$findRange.Find.Execute("misurabile") | Out-Null
$bmRange = $d.Range($findRange.End, $findRange.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Mark the two picture runs (InlineShapes 4 and 5) as NoProofing -> adds <w:rPr><w:noProof/></w:rPr>
$d.InlineShapes(4).Range.NoProofing = 1
$d.InlineShapes(5).Range.NoProofing = 1

Write-Output "done"
